# Update benchmark: 2025-10-11 06:34:42 UTC
#
# Clears the "AKBANK" (column D) benchmark figures for several rows (they no
# longer have a confirmed value) and refreshes the "ISBANKASI" (column E) and
# "FINASNBANK" (column K) incoming-SWIFT maximum fee figures on row 13.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D: clear stale benchmark values, keeping the existing cell style.
$ws.Range("D3").ClearContents()
$ws.Range("D4").ClearContents()
$ws.Range("D5").ClearContents()
$ws.Range("D6").ClearContents()
$ws.Range("D8").ClearContents()
$ws.Range("D9").ClearContents()
$ws.Range("D10").ClearContents()
$ws.Range("D11").ClearContents()
$ws.Range("D12").ClearContents()
$ws.Range("D13").ClearContents()
$ws.Range("D14").ClearContents()

# Row 13: update the refreshed "Azami" (maximum) figures.
$ws.Range("E13").Value = "Hesaba: Asgari 1 TL | Azami 851,5 TL"
$ws.Range("K13").Value = "Hesaba: Asgari 1 TL | Azami 865,75 TL"
